$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'26.236.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -3.96%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'1.659.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.18%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'1.011"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.43%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'217.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.38%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'0.5144"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.60%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'1.010"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.40%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.2603"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.53%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'0.06466"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -3.16%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'19.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -4.60%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.07847"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.46%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("B12").Value = "'Polkadot"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'4.305"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -4.06%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("B13").Value = "'WrappedEther"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'1.650.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.95%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'1.880.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.58%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'0.5532"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.38%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'0.0₅8020"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.91%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'64.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -5.07%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'26.231.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.06%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'1.011"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.40%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'209.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -4.32%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'4.419"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -4.87%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'10.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.80%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'6.048"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.33%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'1.010"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.37%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("B25").Value = "'Toncoin"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'1.811"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +6.08%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("B26").Value = "'Monero"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'144.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.51%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'0.1176"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.63%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'7.003"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.20%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'15.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.83%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'0.05090"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -5.42%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'1.246"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.85%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'3.361"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.28%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'3.249"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -4.12%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'1.561"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.53%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'2.744"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.63%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'2.361"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.56%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'0.9248"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.88%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'1.172.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.01%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.5712"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.80%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.01591"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.59%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("B41").Value = "'mCoin"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'2.568"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.36%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("B42").Value = "'PaxDollar"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'1.011"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.45%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'5.669"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.86%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.8290"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.17%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'100.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.53%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'1.793.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.40%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("E47").Value = "'  -5.97%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'0.4556"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.08%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'55.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.57%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'1.006"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.08%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'7.885"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.93%  "
$ws.Range("E51").Style = "Normal"
